# edit.ps1 - apply the tracked changes to the presentation
#
# Summary of the changes (see commit "Add files via upload"):
#  1. Slide 11 ("Сравнение исходного и оптимизированного дизайна"):
#       - resize/reposition the two comparison screenshots.
#  2. Slide 9 ("Модули программы генетического алгоритма"):
#       - resize the diagram picture (offset stays the same).
#  3. Slide 2 ("Цели и задачи работы"):
#       - merge the "Реализация генетического алгоритма ..." bullet with the
#         following "Реализация алгоритма на языке программирования Python"
#         bullet into a single numbered item.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top/Width/Height are automation `Single`
# (32-bit float) values expressed in points (1 pt = 12700 EMU). Converting an
# EMU value to points and back can land one EMU short because of the float32
# round-trip, so nudge the point value by a hair so it still lands exactly on
# the desired EMU amount once PowerPoint converts it back.
# ---------------------------------------------------------------------------
function EmuToPt([double]$emu) {
    $pts = $emu / 12700.0
    for ($i = 0; $i -lt 50000; $i++) {
        $f32 = [single]$pts
        $back = [math]::Floor([double]$f32 * 12700.0)
        if ($back -eq $emu) {
            return $pts
        } elseif ($back -lt $emu) {
            $pts += 0.0000001
        } else {
            $pts -= 0.0000001
        }
    }
    return $pts
}

# Helper: find a shape on a slide by its (stable) nvSpPr/nvPicPr Id, falling
# back to a 1-based collection index if no shape has that id.
function GetShapeById($slide, [int]$id, [int]$fallbackIndex) {
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $sh = $slide.Shapes.Item($j)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $slide.Shapes.Item($fallbackIndex)
}

# ---------------------------------------------------------------------------
# 1. Slide 11 - reposition/resize the two pictures
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)

$pic1 = GetShapeById $s11 4 2   # "Объект 3"
$pic1.Left   = EmuToPt 395726
$pic1.Top    = EmuToPt 1716505
$pic1.Width  = EmuToPt 5662861
$pic1.Height = EmuToPt 3673642

$pic2 = GetShapeById $s11 5 3   # "Рисунок 4"
$pic2.Left   = EmuToPt 6096000
$pic2.Top    = EmuToPt 2066545
$pic2.Width  = EmuToPt 5946648
$pic2.Height = EmuToPt 3072384

# ---------------------------------------------------------------------------
# 2. Slide 9 - resize the diagram picture
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$pic9 = GetShapeById $s9 12 3   # "Рисунок 11"
$pic9.Left   = EmuToPt 2989632
$pic9.Top    = EmuToPt 1373220
$pic9.Width  = EmuToPt 6539379
$pic9.Height = EmuToPt 5245126

# ---------------------------------------------------------------------------
# 3. Slide 2 - merge the two bullet paragraphs into one
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body = GetShapeById $s2 3 2   # "Объект 2" (the bulleted body placeholder)
$tr = $body.TextFrame.TextRange

# Paragraph 6 ("Реализация алгоритма на языке программирования " + "Python")
# -> rewrite its first run so the merged bullet reads:
#    "Реализация генетического алгоритма для автоматической генерации и
#     оптимизации CSS-стилей на языке программирования " + "Python"
$para6 = $tr.Paragraphs(6, 1)
$run1 = $para6.Runs(1, 1)
$oldRun1Text = $run1.Text
$newRun1Text = "Реализация генетического алгоритма для автоматической генерации и оптимизации CSS-стилей на языке программирования "
$run1Range = $tr.Characters($para6.Start, $oldRun1Text.Length)
$run1Range.Text = $newRun1Text

# Paragraph 5 ("Реализация генетического алгоритма для автоматической
# генерации и оптимизации CSS-стилей") is now redundant - delete it outright,
# which merges what used to be paragraph 6 up into its slot.
$tr2 = $body.TextFrame.TextRange
$para5 = $tr2.Paragraphs(5, 1)
$para5.Delete()
